# Generate Report for Handoff
# Updates the localization-status workbook to reflect that b.md has been
# handed off again (new xliff generated) and is no longer in sync, because
# the previous handback is stale relative to the latest source.

$wb = $excel.ActiveWorkbook

$newStatus       = "Ready for handoff"
$overviewDate    = "2016-08-19 08:36:12"

$zhHandoffFile   = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhHandoffDate   = "2016-08-19 08:36:05"

$deHandoffFile   = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$deHandoffDate   = "2016-08-19 08:36:12"

$errorDetail     = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/59eabea66a64117c5b95a30ad1007913ed7d0901/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8ce9f2d2c9c8cbd2c3d2321ecf93273ae8fc01bd/e2e/b.md."

# ---------------------------------------------------------------------
# Overview sheet: row 3 is the b.md entry.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus
$wsOverview.Range("G3").Value = $overviewDate

# ---------------------------------------------------------------------
# zh-cn sheet: row 3 is the b.md entry.
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $newStatus
$wsZhCn.Range("G3").Value = $zhHandoffFile
$wsZhCn.Range("H3").Value = $zhHandoffDate
$wsZhCn.Range("P3").Value = $errorDetail
# 39.17 characters is what Excel/COM reports back for a column whose
# underlying OOXML <col width> is the raw value 40 (matches column G/J in
# this sheet, which already store width="40" verbatim).
$wsZhCn.Columns.Item(16).ColumnWidth = 39.17

# ---------------------------------------------------------------------
# de-de sheet: row 3 is the b.md entry.
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $newStatus
$wsDeDe.Range("G3").Value = $deHandoffFile
$wsDeDe.Range("H3").Value = $deHandoffDate
$wsDeDe.Range("P3").Value = $errorDetail
$wsDeDe.Columns.Item(16).ColumnWidth = 39.17
